# Add a new column R (year 2021 / value 3.6) to the right of the existing
# Q column on sheet1, matching the style of the Q column cells, and move
# the active selection to O9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data point for 2021 in row 4 (years) and row 5 (percentages),
# mirroring the formatting already applied to column Q (copy the cell
# formatting, then set the new value).
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("R4").Value = 2021

$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("R5").Value = 3.6

# Move the selection as recorded in the saved view state.
$ws.Range("O9").Select()
